$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 312, pushing existing rows 312..387 down to 314..389.
$ws.Rows.Item(312).Insert()
$ws.Rows.Item(312).Insert()

# New row 312: Femacal de La Calera / Repollo / Crespo record / "Primera" record dated 2021-11-08
$ws.Cells.Item(312,1).Value = 3
$ws.Cells.Item(312,2).Value = "Femacal de La Calera"
$ws.Cells.Item(312,3).Value = "Coquimbo"
$ws.Cells.Item(312,4).Value = "11/8/2021"
$ws.Cells.Item(312,5).Value = 5
$ws.Cells.Item(312,6).Value = 100112006
$ws.Cells.Item(312,7).Value = "Repollo"
$ws.Cells.Item(312,8).Value = "Crespo record"
$ws.Cells.Item(312,9).Value = "Primera"
$ws.Cells.Item(312,10).Value = 1800
$ws.Cells.Item(312,11).Value = 600
$ws.Cells.Item(312,12).Value = 650
$ws.Cells.Item(312,13).Value = 626
$ws.Cells.Item(312,14).Value = "$/unidad"
$ws.Cells.Item(312,15).Value = "Provincia de Quillota"
$ws.Cells.Item(312,16).Value = 626
$ws.Cells.Item(312,17).Value = 1
$ws.Cells.Item(312,18).Value = "Hortaliza"

# New row 313: Femacal de La Calera / Repollo / Crespo record / "Segunda" record dated 2021-11-08
$ws.Cells.Item(313,1).Value = 3
$ws.Cells.Item(313,2).Value = "Femacal de La Calera"
$ws.Cells.Item(313,3).Value = "Coquimbo"
$ws.Cells.Item(313,4).Value = "11/8/2021"
$ws.Cells.Item(313,5).Value = 5
$ws.Cells.Item(313,6).Value = 100112006
$ws.Cells.Item(313,7).Value = "Repollo"
$ws.Cells.Item(313,8).Value = "Crespo record"
$ws.Cells.Item(313,9).Value = "Segunda"
$ws.Cells.Item(313,10).Value = 900
$ws.Cells.Item(313,11).Value = 500
$ws.Cells.Item(313,12).Value = 500
$ws.Cells.Item(313,13).Value = 500
$ws.Cells.Item(313,14).Value = "$/unidad"
$ws.Cells.Item(313,15).Value = "Provincia de Quillota"
$ws.Cells.Item(313,16).Value = 500
$ws.Cells.Item(313,17).Value = 1
$ws.Cells.Item(313,18).Value = "Hortaliza"
